$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Se eliminan los periodos de mora anteriores y se agregan los nuevos:
# la fila 16 y la fila 18 intercambian su "Periodo Mora" y "Valor Mora"
# (fila 17 queda igual).
$ws.Range("E16").Value = "1902"
$ws.Range("F16").Value = 28124

$ws.Range("E18").Value = "1809"
$ws.Range("F18").Value = 31249
